$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
# Shared by Overview!E2:F2 / E3:F3, zh-cn!C2:C3, de-de!C2:C3
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C3").Value = "Ready for handoff"

# --- Latest Handback DateTime (de-de H column): 2016-10-17 17:10:53 -> 2016-10-17 17:13:45 ---
$wsDeDe.Range("H2").Value = "2016-10-17 17:13:45"
$wsDeDe.Range("H3").Value = "2016-10-17 17:13:45"

# --- Priority: "ht" -> "mt" ---
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# --- Latest Handoff Datetime (zh-cn H column): 2016-10-17 17:10:28 -> 2016-10-17 17:13:22 ---
$wsZhCn.Range("H2").Value = "2016-10-17 17:13:22"
$wsZhCn.Range("H3").Value = "2016-10-17 17:13:22"

# --- Error Detail for the fe8c1d8f file (row 3) in zh-cn and de-de ---
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e96ad0bd66c8e6deaca9616a154516495bf85bce/e2e/fe8c1d8f-4893-432b-9487-0dc66876f48b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7fe3c855ef923578fced399ffe77aa4da66aeaa4/e2e/fe8c1d8f-4893-432b-9487-0dc66876f48b.md."
$wsZhCn.Range("P3").Value = $errorDetail
$wsDeDe.Range("P3").Value = $errorDetail

# --- Column width adjustments ---
$wsOverview.Range("E1").ColumnWidth = 17.2159881591797
$wsOverview.Range("F1").ColumnWidth = 17.2159881591797

$wsZhCn.Range("C1").ColumnWidth = 17.2159881591797
$wsZhCn.Range("P1").ColumnWidth = 40

$wsDeDe.Range("C1").ColumnWidth = 17.2159881591797
$wsDeDe.Range("P1").ColumnWidth = 40
